$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# The 4th paragraph ("Увеличаване на контраста (contrast enhancement)") is
# removed entirely, merging the following paragraph up into its place.
$para4 = $tr.Paragraphs(4, 1)
$para4.Delete()

# The paragraph that used to be 5th ("Намаляване на шума (denoising)") is now
# the 4th paragraph. Split its leading run "Намаляване на шума " into two
# runs: "Намаляване " and "на шума " (same run formatting for both).
$para4b = $tr.Paragraphs(4, 1)
$firstPart = $tr.Characters($para4b.Start, 11)
$firstPart.Text = "Намаляване "
